$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text value that would otherwise be auto-coerced to a number/percent by Excel,
# forcing it to stay a plain text string with the same General-format data style (no quote-prefix).
$scratch = $ws.Range("Z100")
$ws.Range("H2").Copy()
$scratch.PasteSpecial(-4122)

function Set-TextValue([string]$cellRef, [string]$value) {
    $ws.Range($cellRef).Value = "'" + $value
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

$ws.Range("E2").Value = "2026-02-13 19:18:34"
$ws.Range("M2").Value = "1.2 °C 18:46 TU"
$ws.Range("E3").Value = "2026-02-13 19:18:37"
$ws.Range("G3").Value = "179 cm"
Set-TextValue "H3" "87%"
$ws.Range("I3").Value = "5.8 mm"
$ws.Range("E4").Value = "2026-02-13 19:18:39"
$ws.Range("I4").Value = "6.0 mm"
$ws.Range("J4").Value = "994.7 hPa"
$ws.Range("E5").Value = "2026-02-13 19:18:42"
$ws.Range("I5").Value = "0.8 mm"
$ws.Range("E6").Value = "2026-02-13 19:18:45"
Set-TextValue "H6" "76%"
$ws.Range("I6").Value = "5.0 mm"
$ws.Range("J6").Value = "994.7 hPa"
$ws.Range("E7").Value = "2026-02-13 19:18:47"
$ws.Range("I7").Value = "20.0 mm"
$ws.Range("J7").Value = "995.0 hPa"
$ws.Range("L7").Value = "38.9 km/h - 290º 18:39 TU"
$ws.Range("N7").Value = "10.3 °C 18:37 TU"
$ws.Range("O7").Value = "12.9 °C"
$ws.Range("E8").Value = "2026-02-13 19:18:49"
$ws.Range("I8").Value = "20.6 mm"
$ws.Range("J8").Value = "994.9 hPa"
$ws.Range("N8").Value = "6.8 °C 18:55 TU"
$ws.Range("E9").Value = "2026-02-13 19:18:52"
Set-TextValue "H9" "76%"
$ws.Range("I9").Value = "3.8 mm"
$ws.Range("L9").Value = "30.2 km/h - 330º 18:40 TU"
$ws.Range("O9").Value = "9.6 °C"
$ws.Range("E10").Value = "2026-02-13 19:18:54"
Set-TextValue "H10" "87%"
$ws.Range("I10").Value = "18.3 mm"
$ws.Range("O10").Value = "8.9 °C"
$ws.Range("E11").Value = "2026-02-13 19:18:57"
$ws.Range("O11").Value = "1.8 °C"
$ws.Range("E12").Value = "2026-02-13 19:18:58"
$ws.Range("I12").Value = "5.4 mm"
$ws.Range("E13").Value = "2026-02-13 19:18:59"
$ws.Range("E14").Value = "2026-02-13 19:19:00"
Set-TextValue "H14" "84%"
$ws.Range("I14").Value = "19.9 mm"
$ws.Range("E15").Value = "2026-02-13 19:19:01"
$ws.Range("I15").Value = "3.7 mm"
$ws.Range("E16").Value = "2026-02-13 19:19:02"
$ws.Range("I16").Value = "12.8 mm"
$ws.Range("E17").Value = "2026-02-13 19:19:04"
$ws.Range("E18").Value = "2026-02-13 19:19:05"
$ws.Range("I18").Value = "9.4 mm"
$ws.Range("J18").Value = "994.9 hPa"
$ws.Range("E19").Value = "2026-02-13 19:19:06"
$ws.Range("E20").Value = "2026-02-13 19:19:07"
$ws.Range("I20").Value = "21.8 mm"
$ws.Range("E21").Value = "2026-02-13 19:19:08"
$ws.Range("J21").Value = "997.7 hPa"
$ws.Range("E22").Value = "2026-02-13 19:19:09"
$ws.Range("E23").Value = "2026-02-13 19:19:12"
$ws.Range("I23").Value = "9.4 mm"
$ws.Range("L23").Value = "76.0 km/h - 285º 18:58 TU"
$ws.Range("E24").Value = "2026-02-13 19:19:15"
$ws.Range("J24").Value = "995.6 hPa"
$ws.Range("E25").Value = "2026-02-13 19:19:17"
$ws.Range("G25").Value = "113 cm"
Set-TextValue "H25" "81%"
$ws.Range("I25").Value = "8.9 mm"
$ws.Range("E26").Value = "2026-02-13 19:19:20"
$ws.Range("E27").Value = "2026-02-13 19:19:22"
Set-TextValue "H27" "83%"
$ws.Range("E28").Value = "2026-02-13 19:19:25"
Set-TextValue "H28" "79%"
$ws.Range("J28").Value = "995.2 hPa"
$ws.Range("E29").Value = "2026-02-13 19:19:27"
$ws.Range("O29").Value = "11.1 °C"
$ws.Range("E30").Value = "2026-02-13 19:19:29"
Set-TextValue "H30" "78%"
$ws.Range("I30").Value = "5.1 mm"
$ws.Range("J30").Value = "994.6 hPa"
$ws.Range("L30").Value = "38.5 km/h - 25º 18:48 TU"
$ws.Range("E31").Value = "2026-02-13 19:19:32"
$ws.Range("I31").Value = "3.5 mm"
$ws.Range("J31").Value = "993.6 hPa"
$ws.Range("E32").Value = "2026-02-13 19:19:34"
Set-TextValue "H32" "90%"
$ws.Range("E33").Value = "2026-02-13 19:19:37"
$ws.Range("J33").Value = "996.7 hPa"
$ws.Range("E34").Value = "2026-02-13 19:19:40"
$ws.Range("G34").Value = "109 cm"
$ws.Range("I34").Value = "10.6 mm"
$ws.Range("L34").Value = "62.6 km/h - 11º 18:36 TU"
$ws.Range("E35").Value = "2026-02-13 19:19:42"
$ws.Range("I35").Value = "8.0 mm"
$ws.Range("J35").Value = "995.6 hPa"
$ws.Range("N35").Value = "3.6 °C 18:37 TU"
$ws.Range("E36").Value = "2026-02-13 19:19:45"
Set-TextValue "H36" "78%"
$ws.Range("I36").Value = "8.2 mm"
$ws.Range("J36").Value = "994.8 hPa"
$ws.Range("O36").Value = "10.6 °C"
$ws.Range("E37").Value = "2026-02-13 19:19:47"
Set-TextValue "H37" "84%"
$ws.Range("I37").Value = "12.0 mm"
$ws.Range("J37").Value = "996.6 hPa"
$ws.Range("E38").Value = "2026-02-13 19:19:50"
Set-TextValue "H38" "78%"
$ws.Range("I38").Value = "13.5 mm"
$ws.Range("O38").Value = "9.7 °C"
$ws.Range("E39").Value = "2026-02-13 19:19:52"
$ws.Range("I39").Value = "19.1 mm"
$ws.Range("E40").Value = "2026-02-13 19:19:55"
$ws.Range("J40").Value = "998.2 hPa"
$ws.Range("E41").Value = "2026-02-13 19:19:57"
$ws.Range("J41").Value = "995.1 hPa"
$ws.Range("E42").Value = "2026-02-13 19:20:00"
Set-TextValue "H42" "85%"
$ws.Range("I42").Value = "9.5 mm"
$ws.Range("O42").Value = "11.0 °C"
$ws.Range("E43").Value = "2026-02-13 19:20:02"
$ws.Range("E44").Value = "2026-02-13 19:20:05"
$ws.Range("I44").Value = "6.4 mm"
$ws.Range("E45").Value = "2026-02-13 19:20:08"
Set-TextValue "H45" "63%"
$ws.Range("J45").Value = "993.5 hPa"
$ws.Range("E46").Value = "2026-02-13 19:20:10"
$ws.Range("J46").Value = "995.7 hPa"

# cleanup scratch cell (format + content) so it does not appear in the saved workbook
$scratch.Clear()
